$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 2-7 (the MEI WAN BUILDING entries), causing subsequent rows to shift up.
$ws.Rows("2:7").Delete()
